# Add a new "Status Update 22" sheet (a copy of "Status Update 21"), update its
# figures/commentary for the new reporting period, and move the active-tab
# selection from sheet 21 onto the newly created sheet 22.

$wb = $excel.ActiveWorkbook

# --- 1. Clone the last status-update sheet to create this week's sheet -----
$source = $wb.Worksheets.Item("Status Update 21")

# Record + clear the source sheet's own selection/active-tab state so we can
# restore the (now non-active) look the diff expects once the new sheet is
# in place.
$source.Select()
$source.Range("C13").Select()

$source.Copy($null, $source)
$new = $wb.Worksheets.Item("Status Update 21 (2)")
$new.Name = "Status Update 22"

# --- 2. Update the hours actually logged this period -----------------------
$new.Range("J7").Value = 63
$new.Range("K7").Value = 63
# K10 ( =SUM(K3:K9) ) recalculates automatically to 130.5

# --- 3. Update the narrative cells for this status update -------------------
$new.Range("C11").Value = "Worked on learning how to organize the project to make coding easier. Looked at shadows a bit more, and controls."
$new.Range("C12").Value = "Still haven't found a house to move into after the semester is over."
$new.Range("C13").Value = "On schedule"
$new.Range("E13").Value = "Will change my schedule to reflect this this week."
$new.Range("C14").Value = "Actual hours since last update: 5 hours"
$new.Range("C16").Value = "Need to finish changing my schedule to reflect my actual timeline."

# Scripture cell: plain text with the single word "organized" in bold.
$scriptureCell = $new.Range("C17")
$prefix = "Abraham 4:14 - And [James] "
$bolded = "organized"
$suffix = " the lights in the expanse of the [virtual] heaven, and caused them to divide the day from the night;"
$scriptureCell.Value = $prefix + $bolded + $suffix
$scriptureCell.Characters($prefix.Length + 1, $bolded.Length).Font.Bold = $true

# --- 4. Fix up view/selection state -----------------------------------------
# Sheet 21 is no longer the active tab, but keeps its own last selection.
$source.Select()
$source.Range("C13").Select()

# New sheet becomes the active tab with its own selection.
$new.Select()
$new.Range("E20").Select()
